# Update NATMI LR-pair edge statistics (G:T) for rows 2-16 with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 90.85644533333333
$ws.Range("H2").Value = 272.569336
$ws.Range("I2").Value = 0.6604153190201066
$ws.Range("J2").Value = 0.6604153190201065
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.342152
$ws.Range("N2").Value = 10.026456
$ws.Range("O2").Value = 0.7998676999536424
$ws.Range("P2").Value = 0.7998676999536425
$ws.Range("Q2").Value = 303.6560504836906
$ws.Range("R2").Value = 2732.904454353215
$ws.Range("S2").Value = 0.5282448822387636
$ws.Range("T2").Value = 0.5282448822387636

# Row 3
$ws.Range("G3").Value = 90.85644533333333
$ws.Range("H3").Value = 272.569336
$ws.Range("I3").Value = 0.6604153190201066
$ws.Range("J3").Value = 0.6604153190201065
$ws.Range("M3").Value = 0.4565946666666666
$ws.Range("O3").Value = 0.1092754984925182
$ws.Range("P3").Value = 0.1092754984925182
$ws.Range("Q3").Value = 41.48456837149155
$ws.Range("R3").Value = 373.3611153434239
$ws.Range("S3").Value = 0.07216721319801758
$ws.Range("T3").Value = 0.07216721319801758

# Row 4
$ws.Range("G4").Value = 90.85644533333333
$ws.Range("H4").Value = 272.569336
$ws.Range("I4").Value = 0.6604153190201066
$ws.Range("J4").Value = 0.6604153190201065
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3796343333333334
$ws.Range("N4").Value = 1.138903
$ws.Range("O4").Value = 0.09085680155383948
$ws.Range("P4").Value = 0.09085680155383949
$ws.Range("Q4").Value = 34.49222605315644
$ws.Range("R4").Value = 310.4300344784079
$ws.Range("S4").Value = 0.06000322358332542
$ws.Range("T4").Value = 0.06000322358332542

# Row 5
$ws.Range("G5").Value = 38.83541766666666
$ws.Range("I5").Value = 0.282286024436851
$ws.Range("J5").Value = 0.282286024436851
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.342152
$ws.Range("N5").Value = 10.026456
$ws.Range("O5").Value = 0.7998676999536424
$ws.Range("P5").Value = 0.7998676999536425
$ws.Range("Q5").Value = 129.7938688254853
$ws.Range("R5").Value = 1168.144819429368
$ws.Range("S5").Value = 0.2257914730953617
$ws.Range("T5").Value = 0.2257914730953617

# Row 6
$ws.Range("G6").Value = 38.83541766666666
$ws.Range("I6").Value = 0.282286024436851
$ws.Range("J6").Value = 0.282286024436851
$ws.Range("M6").Value = 0.4565946666666666
$ws.Range("O6").Value = 0.1092754984925182
$ws.Range("P6").Value = 0.1092754984925182
$ws.Range("Q6").Value = 17.73204458437244
$ws.Range("S6").Value = 0.03084694603780807
$ws.Range("T6").Value = 0.03084694603780807

# Row 7
$ws.Range("G7").Value = 38.83541766666666
$ws.Range("I7").Value = 0.282286024436851
$ws.Range("J7").Value = 0.282286024436851
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3796343333333334
$ws.Range("N7").Value = 1.138903
$ws.Range("O7").Value = 0.09085680155383948
$ws.Range("P7").Value = 0.09085680155383949
$ws.Range("Q7").Value = 14.74325789560656
$ws.Range("R7").Value = 132.689321060459
$ws.Range("S7").Value = 0.02564760530368125
$ws.Range("T7").Value = 0.02564760530368126

# Row 8
$ws.Range("G8").Value = 5.478069333333333
$ws.Range("H8").Value = 16.434208
$ws.Range("I8").Value = 0.03981886913046884
$ws.Range("J8").Value = 0.03981886913046883
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.342152
$ws.Range("N8").Value = 10.026456
$ws.Range("O8").Value = 0.7998676999536424
$ws.Range("P8").Value = 0.7998676999536425
$ws.Range("Q8").Value = 18.30854037853867
$ws.Range("R8").Value = 164.776863406848
$ws.Range("S8").Value = 0.0318498272661432
$ws.Range("T8").Value = 0.0318498272661432

# Row 9
$ws.Range("G9").Value = 5.478069333333333
$ws.Range("H9").Value = 16.434208
$ws.Range("I9").Value = 0.03981886913046884
$ws.Range("J9").Value = 0.03981886913046883
$ws.Range("M9").Value = 0.4565946666666666
$ws.Range("O9").Value = 0.1092754984925182
$ws.Range("P9").Value = 0.1092754984925182
$ws.Range("Q9").Value = 2.501257241230222
$ws.Range("R9").Value = 22.511315171072
$ws.Range("S9").Value = 0.004351226773640327
$ws.Range("T9").Value = 0.004351226773640326

# Row 10
$ws.Range("G10").Value = 5.478069333333333
$ws.Range("H10").Value = 16.434208
$ws.Range("I10").Value = 0.03981886913046884
$ws.Range("J10").Value = 0.03981886913046883
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3796343333333334
$ws.Range("N10").Value = 1.138903
$ws.Range("O10").Value = 0.09085680155383948
$ws.Range("P10").Value = 0.09085680155383949
$ws.Range("Q10").Value = 2.079663199313778
$ws.Range("R10").Value = 18.716968793824
$ws.Range("S10").Value = 0.003617815090685312
$ws.Range("T10").Value = 0.003617815090685312

# Row 11
$ws.Range("G11").Value = 1.120690333333333
$ws.Range("H11").Value = 3.362071
$ws.Range("I11").Value = 0.008146049092012497
$ws.Range("J11").Value = 0.008146049092012496
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.342152
$ws.Range("N11").Value = 10.026456
$ws.Range("O11").Value = 0.7998676999536424
$ws.Range("P11").Value = 0.7998676999536425
$ws.Range("Q11").Value = 3.745517438930667
$ws.Range("R11").Value = 33.70965695037599
$ws.Range("S11").Value = 0.006515761550937493
$ws.Range("T11").Value = 0.006515761550937493

# Row 12
$ws.Range("G12").Value = 1.120690333333333
$ws.Range("H12").Value = 3.362071
$ws.Range("I12").Value = 0.008146049092012497
$ws.Range("J12").Value = 0.008146049092012496
$ws.Range("M12").Value = 0.4565946666666666
$ws.Range("O12").Value = 0.1092754984925182
$ws.Range("P12").Value = 0.1092754984925182
$ws.Range("Q12").Value = 0.5117012291848889
$ws.Range("R12").Value = 4.605311062664
$ws.Range("S12").Value = 0.0008901635752741908
$ws.Range("T12").Value = 0.0008901635752741907

# Row 13
$ws.Range("G13").Value = 1.120690333333333
$ws.Range("H13").Value = 3.362071
$ws.Range("I13").Value = 0.008146049092012497
$ws.Range("J13").Value = 0.008146049092012496
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3796343333333334
$ws.Range("N13").Value = 1.138903
$ws.Range("O13").Value = 0.09085680155383948
$ws.Range("P13").Value = 0.09085680155383949
$ws.Range("Q13").Value = 0.4254525275681111
$ws.Range("R13").Value = 3.829072748113
$ws.Range("S13").Value = 0.0007401239658008138
$ws.Range("T13").Value = 0.0007401239658008136

# Row 14
$ws.Range("G14").Value = 1.284086333333333
$ws.Range("H14").Value = 3.852259
$ws.Range("I14").Value = 0.009333738320561039
$ws.Range("J14").Value = 0.009333738320561037
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.342152
$ws.Range("N14").Value = 10.026456
$ws.Range("O14").Value = 0.7998676999536424
$ws.Range("P14").Value = 0.7998676999536425
$ws.Range("Q14").Value = 4.291611707122667
$ws.Range("R14").Value = 38.624505364104
$ws.Range("S14").Value = 0.007465755802436331
$ws.Range("T14").Value = 0.00746575580243633

# Row 15
$ws.Range("G15").Value = 1.284086333333333
$ws.Range("H15").Value = 3.852259
$ws.Range("I15").Value = 0.009333738320561039
$ws.Range("J15").Value = 0.009333738320561037
$ws.Range("M15").Value = 0.4565946666666666
$ws.Range("O15").Value = 0.1092754984925182
$ws.Range("P15").Value = 0.1092754984925182
$ws.Range("Q15").Value = 0.5863069713395556
$ws.Range("R15").Value = 5.276762742055999
$ws.Range("S15").Value = 0.001019948907778027
$ws.Range("T15").Value = 0.001019948907778027

# Row 16
$ws.Range("G16").Value = 1.284086333333333
$ws.Range("H16").Value = 3.852259
$ws.Range("I16").Value = 0.009333738320561039
$ws.Range("J16").Value = 0.009333738320561037
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.3796343333333334
$ws.Range("N16").Value = 1.138903
$ws.Range("O16").Value = 0.09085680155383948
$ws.Range("P16").Value = 0.09085680155383949
$ws.Range("Q16").Value = 0.4874832590974445
$ws.Range("R16").Value = 4.387349331877
$ws.Range("S16").Value = 0.0008480336103466813
$ws.Range("T16").Value = 0.0008480336103466813

